$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    # Writes $value into $cellRef as literal text, even if it
    # looks like a number (e.g. "240.32"), by round-tripping
    # through a text NumberFormat and restoring the original style.
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '96.927.94'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '3.701.48'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue $ws 'D5' '240.32'
$ws.Range('E5').Value = '  -0.99%  '
Set-TextValue $ws 'D6' '1.90'
$ws.Range('E6').Value = '  +9.09%  '
Set-TextValue $ws 'D7' '655.26'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('E9').Value = '  +3.90%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').Value = '3.699.96'
$ws.Range('E11').Value = '  +2.95%  '
Set-TextValue $ws 'D12' '45.42'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('E13').Value = '  +1.04%  '
Set-TextValue $ws 'D14' '6.88'
$ws.Range('E14').Value = '  +6.06%  '
$ws.Range('D15').Value = '4.388.25'
$ws.Range('E15').Value = '  +2.82%  '
Set-TextValue $ws 'D16' '0.0000269'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '96.659.40'
$ws.Range('E17').Value = '  -0.63%  '
Set-TextValue $ws 'D18' '9.09'
$ws.Range('E18').Value = '  +4.70%  '
$ws.Range('D19').Value = '3.684.29'
$ws.Range('E19').Value = '  +2.45%  '
Set-TextValue $ws 'D20' '19.30'
$ws.Range('E20').Value = '  +6.54%  '
Set-TextValue $ws 'D21' '12.92'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('E22').Value = '  -0.26%  '
Set-TextValue $ws 'D23' '526.93'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('E24').Value = '  +0.59%  '
Set-TextValue $ws 'D25' '7.13'
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('E26').Value = '  -3.08%  '
Set-TextValue $ws 'D27' '102.00'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('D29').Value = '3.896.28'
$ws.Range('E29').Value = '  +2.78%  '
$ws.Range('E30').Value = '  -1.38%  '
Set-TextValue $ws 'D31' '12.60'
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  +13.97%  '
$ws.Range('E35').Value = '  +0.36%  '
Set-TextValue $ws 'D36' '32.81'
$ws.Range('E36').Value = '  +2.70%  '
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D38' '659.52'
$ws.Range('E38').Value = '  +6.73%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws 'D39' '0.608'
$ws.Range('E39').Value = '  +6.18%  '
Set-TextValue $ws 'D40' '9.05'
$ws.Range('E40').Value = '  +3.62%  '
Set-TextValue $ws 'D41' '7.12'
$ws.Range('E41').Value = '  +18.03%  '
$ws.Range('E42').Value = '  +5.54%  '
$ws.Range('E43').Value = '  +3.49%  '
Set-TextValue $ws 'D44' '0.970'
$ws.Range('E44').Value = '  +4.14%  '
Set-TextValue $ws 'D45' '39.10'
$ws.Range('E45').Value = '  +18.69%  '
$ws.Range('E46').Value = '  +0.05%  '
Set-TextValue $ws 'D47' '0.456'
$ws.Range('E47').Value = '  +5.22%  '
Set-TextValue $ws 'D48' '0.0460'
$ws.Range('E48').Value = '  +4.19%  '
$ws.Range('E49').Value = '  +0.54%  '
Set-TextValue $ws 'D50' '8.83'
$ws.Range('E50').Value = '  +2.90%  '
$ws.Range('E51').Value = '  -0.02%  '
